$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "Added SVR parameter loading from pred_par structure and Excel files"
# -> add three new SVR hyper-parameter columns (header in row 1, value in row 2)
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 0.05
$ws.Range("M2").Value = 20

# Leave the selection on the newly added cell, matching the saved view
$ws.Range("K2").Select()
